$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates scraped from the "Updated symbol list" GitHub Actions commit.
# Columns B (Coin), C (Link), D (Price) and E (Volume 1h) hold text values
# (including numeric-looking strings and percentages), so each target cell
# is formatted as Text before the value is written. This prevents Excel
# from re-interpreting strings like "0.1170" or "34.18%" as numbers, which
# would silently drop trailing zeros or convert percentages into floats.
$updates = @(
    @{ Cell = 'D2'; Value = '275.61' },
    @{ Cell = 'E2'; Value = '-1.22%' },
    @{ Cell = 'D3'; Value = '26.57' },
    @{ Cell = 'E3'; Value = '-2.66%' },
    @{ Cell = 'D4'; Value = '4.882' },
    @{ Cell = 'E4'; Value = '1.80%' },
    @{ Cell = 'D5'; Value = '0.06337' },
    @{ Cell = 'E5'; Value = '0.32%' },
    @{ Cell = 'D6'; Value = '6.912' },
    @{ Cell = 'E6'; Value = '-0.22%' },
    @{ Cell = 'D7'; Value = '3.315' },
    @{ Cell = 'E7'; Value = '0.63%' },
    @{ Cell = 'E8'; Value = '34.18%' },
    @{ Cell = 'D9'; Value = '0.8714' },
    @{ Cell = 'E9'; Value = '-0.63%' },
    @{ Cell = 'D10'; Value = '0.1544' },
    @{ Cell = 'E10'; Value = '5.10%' },
    @{ Cell = 'D11'; Value = '0.05024' },
    @{ Cell = 'E11'; Value = '-1.02%' },
    @{ Cell = 'D12'; Value = '0.07404' },
    @{ Cell = 'E12'; Value = '1.12%' },
    @{ Cell = 'D13'; Value = '0.02958' },
    @{ Cell = 'E13'; Value = '-6.08%' },
    @{ Cell = 'D14'; Value = '0.09045' },
    @{ Cell = 'E14'; Value = '-0.27%' },
    @{ Cell = 'E15'; Value = '1.02%' },
    @{ Cell = 'D16'; Value = '0.0006322' },
    @{ Cell = 'E16'; Value = '0.99%' },
    @{ Cell = 'D17'; Value = '0.006017' },
    @{ Cell = 'E17'; Value = '2.27%' },
    @{ Cell = 'D18'; Value = '3.448' },
    @{ Cell = 'E18'; Value = '0.00%' },
    @{ Cell = 'E19'; Value = '-0.55%' },
    @{ Cell = 'E21'; Value = '1.11%' },
    @{ Cell = 'D22'; Value = '3.902' },
    @{ Cell = 'E22'; Value = '0.89%' },
    @{ Cell = 'D23'; Value = '0.04372' },
    @{ Cell = 'E23'; Value = '0.74%' },
    @{ Cell = 'D24'; Value = '0.001178' },
    @{ Cell = 'E24'; Value = '-0.23%' },
    @{ Cell = 'E25'; Value = '-1.62%' },
    @{ Cell = 'D26'; Value = '0.0001199' },
    @{ Cell = 'E26'; Value = '-0.06%' },
    @{ Cell = 'E40'; Value = '0.77%' },
    @{ Cell = 'D41'; Value = '0.006953' },
    @{ Cell = 'E41'; Value = '5.51%' },
    @{ Cell = 'D42'; Value = '0.1170' },
    @{ Cell = 'E42'; Value = '0.97%' },
    @{ Cell = 'E43'; Value = '-2.54%' },
    @{ Cell = 'E44'; Value = '-12.24%' },
    @{ Cell = 'D45'; Value = '0.00005294' },
    @{ Cell = 'E45'; Value = '1.33%' },
    @{ Cell = 'B46'; Value = 'CoinbaseStockToken' },
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin' },
    @{ Cell = 'D46'; Value = '0.02099' },
    @{ Cell = 'E46'; Value = '-6.72%' },
    @{ Cell = 'B47'; Value = 'BOLO' },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo' },
    @{ Cell = 'D47'; Value = '1.490' },
    @{ Cell = 'E47'; Value = '-37.36%' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
